$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customers 1")

$ws.Range("F2").Value = "0731264413"
$ws.Range("F3").Value = "0731212345"
$ws.Range("F4").Value = "0731234412"
$ws.Range("F5").Value = "0731233332"
$ws.Range("F6").Value = "0761111234"
$ws.Range("F7").Value = "0751131244"
$ws.Range("F8").Value = "08444555"
$ws.Range("F9").Value = "085556666"
$ws.Range("F10").Value = "081231234"
$ws.Range("F11").Value = "085556666"
